$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 49000
$ws.Range("I13").Value = 50000
$ws.Range("J13").Value = 48000
$ws.Range("K13").Value = 50000
$ws.Range("L13").Value = 48000
$ws.Range("M13").Value = -49831
$ws.Range("N13").Value = -48338
$ws.Range("H55").Value = 71767.92999999999
$ws.Range("J55").Value = 356.5
$ws.Range("L55").Value = 356.5
$ws.Range("N55").Value = -784.5
$ws.Range("H86").Value = 1584.7142
$ws.Range("I86").Value = 1432.1666
$ws.Range("J86").Value = 2500
$ws.Range("K86").Value = 1432.1666
$ws.Range("L86").Value = 2500
$ws.Range("M86").Value = -309.1666
$ws.Range("N86").Value = -4746
$ws.Range("H89").Value = 1584.7142
$ws.Range("I89").Value = 1432.1666
$ws.Range("J89").Value = 2500
$ws.Range("K89").Value = 7160.833000000001
$ws.Range("L89").Value = 12500
$ws.Range("M89").Value = -1544.833000000001
$ws.Range("N89").Value = -23732
$ws.Range("H106").Value = 5606841
$ws.Range("I106").Value = 5901817
$ws.Range("J106").Value = 2300
$ws.Range("K106").Value = 5901817
$ws.Range("L106").Value = 2300
$ws.Range("M106").Value = -5901186
$ws.Range("N106").Value = -3562
$ws.Range("H137").Value = 111115840
$ws.Range("I137").Value = 166670830
$ws.Range("J137").Value = 5866.6665
$ws.Range("K137").Value = 500012490
$ws.Range("L137").Value = 17599.9995
$ws.Range("M137").Value = -500009940
$ws.Range("N137").Value = -22699.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 8128.3335
$ws.Range("J43").Value = 8128.3335
$ws.Range("L43").Value = 8128.3335
$ws.Range("N43").Value = -8754.333500000001
$ws.Range("H61").Value = 3693
$ws.Range("I61").Value = 2734.5557
$ws.Range("K61").Value = 2734.5557
$ws.Range("M61").Value = -2522.5557
$ws.Range("H136").Value = 3693
$ws.Range("I136").Value = 2734.5557
$ws.Range("K136").Value = 8203.667099999999
$ws.Range("M136").Value = -5653.667099999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 658.0769
$ws.Range("I94").Value = 623.2727
$ws.Range("J94").Value = 849.5
$ws.Range("K94").Value = 623.2727
$ws.Range("L94").Value = 849.5
$ws.Range("M94").Value = -172.2727
$ws.Range("N94").Value = -1751.5
$ws.Range("H105").Value = 2482.8438
$ws.Range("I105").Value = 2394.2693
$ws.Range("J105").Value = 2866.6667
$ws.Range("K105").Value = 2394.2693
$ws.Range("L105").Value = 2866.6667
$ws.Range("M105").Value = -647.2692999999999
$ws.Range("N105").Value = -6360.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1666.3529
$ws.Range("I58").Value = 1157.8462
$ws.Range("J58").Value = 3319
$ws.Range("K58").Value = 1157.8462
$ws.Range("L58").Value = 3319
$ws.Range("M58").Value = -954.8462
$ws.Range("N58").Value = -3725
$ws.Range("H136").Value = 1666.3529
$ws.Range("I136").Value = 1157.8462
$ws.Range("J136").Value = 3319
$ws.Range("K136").Value = 3473.5386
$ws.Range("L136").Value = 9957
$ws.Range("M136").Value = -923.5385999999999
$ws.Range("N136").Value = -15057

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 287.875
$ws.Range("I107").Value = 232.25
$ws.Range("J107").Value = 315.6875
$ws.Range("K107").Value = 696.75
$ws.Range("L107").Value = 947.0625
$ws.Range("M107").Value = 1223.25
$ws.Range("N107").Value = -4787.0625
$ws.Range("H131").Value = 1423.0299
$ws.Range("I131").Value = 493.9
$ws.Range("J131").Value = 1586.035
$ws.Range("K131").Value = 1481.7
$ws.Range("L131").Value = 4758.105
$ws.Range("M131").Value = 3558.3
$ws.Range("N131").Value = -14838.105
$ws.Range("H134").Value = 2101.2942
$ws.Range("I134").Value = 1482.625
$ws.Range("K134").Value = 4447.875
$ws.Range("M134").Value = 622.125
$ws.Range("H137").Value = 5614565
$ws.Range("I137").Value = 8335842
$ws.Range("J137").Value = 172011
$ws.Range("K137").Value = 25007526
$ws.Range("L137").Value = 516033
$ws.Range("M137").Value = -25002426
$ws.Range("N137").Value = -526233

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 19980
$ws.Range("J69").Value = 19980
$ws.Range("L69").Value = 19980
$ws.Range("N69").Value = -21478
$ws.Range("H72").Value = 19980
$ws.Range("J72").Value = 19980
$ws.Range("L72").Value = 59940
$ws.Range("N72").Value = -67428
$ws.Range("H74").Value = 59798.75
$ws.Range("J74").Value = 59798.75
$ws.Range("L74").Value = 59798.75
$ws.Range("N74").Value = -61670.75
$ws.Range("H75").Value = 29993.334
$ws.Range("J75").Value = 29993.334
$ws.Range("L75").Value = 29993.334
$ws.Range("N75").Value = -31741.334
$ws.Range("H77").Value = 59798.75
$ws.Range("J77").Value = 59798.75
$ws.Range("L77").Value = 179396.25
$ws.Range("N77").Value = -188756.25
$ws.Range("H78").Value = 29993.334
$ws.Range("J78").Value = 29993.334
$ws.Range("L78").Value = 89980.00199999999
$ws.Range("N78").Value = -98716.00199999999
$ws.Range("H86").Value = 34890
$ws.Range("J86").Value = 34890
$ws.Range("L86").Value = 34890
$ws.Range("N86").Value = -37262
$ws.Range("H88").Value = 23445
$ws.Range("J88").Value = 23445
$ws.Range("L88").Value = 23445
$ws.Range("N88").Value = -24347
$ws.Range("H89").Value = 34890
$ws.Range("J89").Value = 34890
$ws.Range("L89").Value = 104670
$ws.Range("N89").Value = -116526
$ws.Range("H91").Value = 23445
$ws.Range("J91").Value = 23445
$ws.Range("L91").Value = 23445
$ws.Range("N91").Value = -26565
$ws.Range("H113").Value = 2051.1
$ws.Range("I113").Value = 2027.75
$ws.Range("J113").Value = 2066.6667
$ws.Range("K113").Value = 2027.75
$ws.Range("L113").Value = 2066.6667
$ws.Range("M113").Value = 142.25
$ws.Range("N113").Value = -6406.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 11822
$ws.Range("I29").Value = 13733
$ws.Range("J29").Value = 8000
$ws.Range("K29").Value = 13733
$ws.Range("L29").Value = 8000
$ws.Range("M29").Value = -13438
$ws.Range("N29").Value = -8590
$ws.Range("H122").Value = 3433.3333
$ws.Range("I122").Value = 1800
$ws.Range("J122").Value = 3529.4119
$ws.Range("K122").Value = 5400
$ws.Range("L122").Value = 10588.2357
$ws.Range("M122").Value = -2950
$ws.Range("N122").Value = -15488.2357

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H123").Value = 30429
$ws.Range("J123").Value = 30429
$ws.Range("L123").Value = 30429
$ws.Range("N123").Value = -40229
$ws.Range("H132").Value = 50004652
$ws.Range("I132").Value = 62504812
$ws.Range("J132").Value = 4002
$ws.Range("K132").Value = 187514436
$ws.Range("L132").Value = 12006
$ws.Range("M132").Value = -187511906
$ws.Range("N132").Value = -17066
